# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets
# to reflect the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 25
$ws1.Range("F4").Value  = 88
$ws1.Range("F5").Value  = 22
$ws1.Range("F7").Value  = 1712
$ws1.Range("F8").Value  = 32
$ws1.Range("F11").Value = 1777
$ws1.Range("F13").Value = 110
$ws1.Range("F14").Value = 423
$ws1.Range("F18").Value = 33
$ws1.Range("F21").Value = 741
$ws1.Range("F22").Value = 307
$ws1.Range("F23").Value = 166
$ws1.Range("F24").Value = 242

# Sheet "全部类型" (sheetId 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 25
$ws4.Range("F4").Value  = 88
$ws4.Range("F5").Value  = 22
$ws4.Range("F7").Value  = 1712
$ws4.Range("F9").Value  = 32
$ws4.Range("F12").Value = 1777
$ws4.Range("F14").Value = 110
$ws4.Range("F15").Value = 423
$ws4.Range("F19").Value = 33
$ws4.Range("F22").Value = 741
$ws4.Range("F23").Value = 307
$ws4.Range("F24").Value = 166
$ws4.Range("F25").Value = 242

$wb.Save()
